$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Nine")

$ws.Range("C2").Value = 10.11

$ws.Range("B3").Value = 9.890000000000001
$ws.Range("D3").Value = 10.44
$ws.Range("E3").Value = 10.74

$ws.Range("C4").Value = 9.56
$ws.Range("E4").Value = 10.67
$ws.Range("F4").Value = 10.19

$ws.Range("C5").Value = 9.26
$ws.Range("D5").Value = 9.33
$ws.Range("G5").Value = 9.75

$ws.Range("D6").Value = 9.81
$ws.Range("G6").Value = 10.32
$ws.Range("H6").Value = 10.47

$ws.Range("E7").Value = 10.25
$ws.Range("F7").Value = 9.68
$ws.Range("H7").Value = 9.81
$ws.Range("I7").Value = 6.3

$ws.Range("F8").Value = 9.529999999999999
$ws.Range("G8").Value = 10.19
$ws.Range("I8").Value = 8.67

$ws.Range("G9").Value = 13.7
$ws.Range("H9").Value = 11.33
